$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: header + width ---
$ws.Range("E1").Value = "Areas of Impact"
$ws.Columns("E").ColumnWidth = 51.95

# --- Column C gets "Codex Audit - Architecture musts" for the existing Pattern rows (19-27) ---
$ws.Range("C19").Value = "Codex Audit - Architecture musts"
$ws.Range("C20").Value = "Codex Audit - Architecture musts"
$ws.Range("C21").Value = "Codex Audit - Architecture musts"
$ws.Range("C22").Value = "Codex Audit - Architecture musts"
$ws.Range("C23").Value = "Codex Audit - Architecture musts"
$ws.Range("C24").Value = "Codex Audit - Architecture musts"
$ws.Range("C25").Value = "Codex Audit - Architecture musts"
$ws.Range("C26").Value = "Codex Audit - Architecture musts"
$ws.Range("C27").Value = "Codex Audit - Architecture musts"

# --- Row 28: new Flow 2 row with a red "To Do" note in D ---
$ws.Range("A28").Value = "## Critical Implementation Patterns"
$ws.Range("B28").Value = "### Flow 2: Automatic Creator Onboarding"
$ws.Range("D28").Value = "Check after signup API Contract is created"
$ws.Range("D28").Font.Color = 255

# --- Row 29: Discovery & Onboarding Model reused as a new data row ---
$ws.Range("A29").Value = "## Critical Implementation Patterns"
$ws.Range("B29").Value = "### Discovery & Onboarding Model"
$ws.Range("C29").Value = "How creators join program. Two Paths: Cruva-Sourced Creators + Word-of-outh Creators"
$ws.Range("E29").Value = "auth pages flow"

# --- Row 30: Flow 3 ---
$ws.Range("A30").Value = "## Critical Implementation Patterns"
$ws.Range("B30").Value = "### Flow 3: Creator First-Time Registration"
$ws.Range("C30").Value = "Creator Flows"
$ws.Range("E30").Value = "auth pages flow"

# --- Row 31: Flow 4 ---
$ws.Range("A31").Value = "## Critical Implementation Patterns"
$ws.Range("B31").Value = "### Flow 4: Returning User Login"
$ws.Range("C31").Value = "Creator Flows"
$ws.Range("E31").Value = "auth pages flow"

# --- Row 32: Flow 5 ---
$ws.Range("A32").Value = "## Critical Implementation Patterns"
$ws.Range("B32").Value = "### Flow 5: Password Reset (Magic Link)"
$ws.Range("C32").Value = "Creator Flows"
$ws.Range("E32").Value = "auth pages flow"

# --- Row 33: Flow 6 ---
$ws.Range("A33").Value = "## Critical Implementation Patterns"
$ws.Range("B33").Value = "### Flow 6: Email Verification System (OTP)"
$ws.Range("C33").Value = "Creator Flows"
$ws.Range("E33").Value = "auth pages flow"

# --- Row 34: Flow 7 (daily tier calc) with a red note in D ---
$ws.Range("A34").Value = "## Critical Implementation Patterns"
$ws.Range("B34").Value = "### Flow 7: Daily Tier Calculation (Automated)"
$ws.Range("D34").Value = "This is a daily flow, automated. Not related to creators login flows"
$ws.Range("D34").Font.Color = 255

# --- Row 35: Flow 7 (admin adds creator manually) ---
$ws.Range("A35").Value = "## Critical Implementation Patterns"
$ws.Range("B35").Value = "### Flow 7: Admin Adds Creator Manually"
$ws.Range("C35").Value = "Creator Flows"
$ws.Range("E35").Value = "auth pages flow"

# --- Row 36: Flow 8 (creator claims reward) ---
$ws.Range("A36").Value = "## Critical Implementation Patterns"
$ws.Range("B36").Value = "### Flow 8: Creator Claims Reward"

# --- View state: freeze header row, scroll down, select C37 ---
$ws.Range("A2").Select()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.Split = $false
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C37").Select()

Write-Host "Edit complete"
